$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 370, shifting existing rows 370-487 down to 371-488
$ws.Rows.Item(370).Insert()

# Populate the newly inserted row 370 with the new data record
$ws.Cells.Item(370, 1).Value = 5
$ws.Cells.Item(370, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(370, 3).Value = "Maule"
$ws.Cells.Item(370, 4).Value = 44988
$ws.Cells.Item(370, 5).Value = 7
$ws.Cells.Item(370, 6).Value = 100114014
$ws.Cells.Item(370, 7).Value = "Betarraga"
$ws.Cells.Item(370, 8).Value = "Sin especificar"
$ws.Cells.Item(370, 9).Value = "Primera"
$ws.Cells.Item(370, 10).Value = 4000
$ws.Cells.Item(370, 11).Value = 600
$ws.Cells.Item(370, 12).Value = 600
$ws.Cells.Item(370, 13).Value = 600
$ws.Cells.Item(370, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(370, 15).Value = "Región del Maule"
$ws.Cells.Item(370, 16).Value = 120
$ws.Cells.Item(370, 17).Value = 5
$ws.Cells.Item(370, 18).Value = "Hortaliza"
